# GEO-RES.schema.docx edits
# Strategy notes:
#  - Plain text replacements (no <w:br/> needed) are done via Cell.Range.Text = "..."
#    which preserves straight apostrophes and auto-adds xml:space="preserve" when needed.
#  - Text that needs a manual line break (<w:br/>) is built with a unique "|BREAK|"
#    placeholder via Cell.Range.Text = "...", then a narrow Find/Replace (scoped through
#    a fresh $d.Range(start,end) built right after the Cell.Range is re-fetched) swaps
#    "|BREAK|" for "^l" (Word's manual line-break code). Using Find.Execute only for this
#    tiny placeholder (never for the real natural-language text) avoids Word's
#    smart-quote autocorrect from curling straight apostrophes in the surrounding text.

$d = $word.ActiveDocument

function Insert-Break {
    param($cellRangeStart, $cellRangeEnd)
    $scoped = $d.Range($cellRangeStart, $cellRangeEnd)
    $scoped.Find.Execute("|BREAK|", $false, $false, $false, $false, $false, $true, 0, $false, "^l", 1) | Out-Null
}

$t2 = $d.Tables(2)
$t3 = $d.Tables(3)

# 1) id / Exemple: 76_45101#SMUR1# -> 76_45101#SMUR1
$t2.Cell(2,6).Range.Text = "76_45101#SMUR1"

# 2) orgId / Description: drop "{code département}." + add "*données facultatives" line
$t2.Cell(3,5).Range.Text = "Identifiant unique de l'organisme :  {pays}.{domaine}.{organisation}.{structure interne}*.{unité fonctionnelle}*|BREAK|*données facultatives"
$c = $t2.Cell(3,5).Range
Insert-Break $c.Start $c.End

# 3) orgId / Exemple: fr.health.760.samu76 -> fr.health.samu76A
$t2.Cell(3,6).Range.Text = "fr.health.samu76A"

# 4) name / Champ correspondant: new two-line description
$t2.Cell(4,2).Range.Text = "Nom donné à la ressource par l'organisme propriétaire. |BREAK|L'immatriculation peut être utilisée dans le nom courant des véhicules."
$c = $t2.Cell(4,2).Range
Insert-Break $c.Start $c.End

# 5) type / Format: string -> string + ENUM list
$t2.Cell(5,3).Range.Text = "string|BREAK|(ENUM: SMUR, SDIS, TSU, SNP, MSPE, SHIP)"
$c = $t2.Cell(5,3).Range
Insert-Break $c.Start $c.End

# 6) nature / Format: string -> string + ENUM list
$t2.Cell(6,3).Range.Text = "string|BREAK|(ENUM: EFFECTEUR, BASE)"
$c = $t2.Cell(6,3).Range
Insert-Break $c.Start $c.End

# 7) nature / Description: effector -> effecteur
$t2.Cell(6,5).Range.Text = "Nature de la ressource (effecteur, base)"

# 8) nature / Exemple: (empty) -> BASE
$t2.Cell(6,6).Range.Text = "BASE"

# 9) mobility / Format: string -> string + ENUM list
$t2.Cell(7,3).Range.Text = "string|BREAK|(ENUM: FIXE, VEHICULE, HELICOPTERE, SHIP )"
$c = $t2.Cell(7,3).Range
Insert-Break $c.Start $c.End

# 10) mobility / Description: translated wording
$t2.Cell(7,5).Range.Text = "Mobilité de la ressource (fixe, vehicule, heliporté, navire)"

# 11) mobility / Exemple: Vehicle -> VEHCL
$t2.Cell(7,6).Range.Text = "VEHCL"

# 12) capacity / Format: string -> string + ENUM list
$t2.Cell(8,3).Range.Text = "string|BREAK|(ENUM: URGENCE, MEDICALE, PARAMEDICALE, INCONNUE)"
$c = $t2.Cell(8,3).Range
Insert-Break $c.Start $c.End

# 13) capacity / Exemple: (empty) -> MED
$t2.Cell(8,6).Range.Text = "MED"

# 14) (Type contact table) type / Format: string -> string + ENUM list
$t3.Cell(2,3).Range.Text = "string|BREAK|(ENUM: PMRADD, PHNADD)"
$c = $t3.Cell(2,3).Range
Insert-Break $c.Start $c.End
